$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Hoja1. Actividades")
$ws2 = $wb.Worksheets.Item("Hoja2. Integrantes")

# --- Hoja1. Actividades ---

# Director name: remove accent on "Jose"
$ws1.Range("E5").Value = "Jose Alejandro Cortés Taborda"

# Row 10: clear the "X" marks and the extra note in J10
$ws1.Range("D10:E10").ClearContents()
$ws1.Range("J10").ClearContents()

# Row 11: clear the whole test/demo row
$ws1.Range("A11:C11").ClearContents()
$ws1.Range("F11:G11").ClearContents()
$ws1.Range("I11:J11").ClearContents()

# --- Hoja2. Integrantes ---

# Row 6: replace first member with the director's info
$ws2.Range("A6").Value = "Jose Alejandro Cortés Taborda"
$ws2.Range("B6").Value = "Director"
$ws2.Range("C6").Value = 1035
$ws2.Range("D6").Value = 12
$ws2.Range("E6").Value = "jose_cortes82141@elpoli.edu.co"

# Row 7: clear the second member entirely
$ws2.Range("A7:E7").ClearContents()
